$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column C ("full"), shifting
# the old C (full), D (tipo), E (link) columns to E, F, G respectively.
$ws.Range("C1:D1").EntireColumn.Insert()

# New header row values
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Row 2 - Fonte 90 Bob
$ws.Range("C2").Value = "FONTE 90 BOB"
$ws.Range("D2").Value = "Igual"
$ws.Range("F2").Value = "premium"
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:442970967#searchVariation=MLB21562641&position=2&search_layout=stack&type=product&tracking_id=3ecc166e-cc7f-401e-8926-b1671e3d6284"

# Row 3 - Fonte 120A
$ws.Range("C3").Value = "FONTE 120A"
$ws.Range("D3").Value = "Igual"
$ws.Range("F3").Value = "premium"
$ws.Range("G3").Value = "https://www.mercadolivre.com.br/fonte-automotiva-120a-amperes-jfa-carregador-cor-preto/p/MLB21392652?pdp_filters=seller_id:442970967#searchVariation=MLB21392652&position=6&search_layout=stack&type=product&tracking_id=3ecc166e-cc7f-401e-8926-b1671e3d6284"

# Row 4 - Sem Modelo (no politica value)
$ws.Range("C4").Value = "Sem Modelo"
$ws.Range("D4").Value = ""
$ws.Range("F4").Value = "classico"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:442970967#searchVariation=MLB27687422&position=7&search_layout=stack&type=product&tracking_id=3ecc166e-cc7f-401e-8926-b1671e3d6284"

# Row 5 - Fonte 200 Bob
$ws.Range("C5").Value = "FONTE 200 BOB"
$ws.Range("D5").Value = "Igual"
$ws.Range("F5").Value = "premium"
$ws.Range("G5").Value = "https://www.mercadolivre.com.br/fonte-automotiva-jfa-storm-200a-bob-carregador-automatico-bivolt-cor-bob-200a-jfa/p/MLB24834408?pdp_filters=seller_id:442970967#searchVariation=MLB24834408&position=5&search_layout=stack&type=product&tracking_id=3ecc166e-cc7f-401e-8926-b1671e3d6284"

# Row 6 - Fonte 120 Bob
$ws.Range("C6").Value = "FONTE 120 BOB"
$ws.Range("D6").Value = "Igual"
$ws.Range("F6").Value = "classico"
$ws.Range("G6").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-120a-bob-slim-bivolt-cor-preto/p/MLB22144397?pdp_filters=seller_id:442970967#searchVariation=MLB22144397&position=3&search_layout=stack&type=product&tracking_id=3ecc166e-cc7f-401e-8926-b1671e3d6284"

# Row 7 - Fonte 200A
$ws.Range("C7").Value = "FONTE 200A"
$ws.Range("D7").Value = "Igual"
$ws.Range("F7").Value = "premium"
$ws.Range("G7").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotiva-jfa-200a-slim-bivolt-voltimetro/p/MLB21348561?pdp_filters=seller_id:442970967#searchVariation=MLB21348561&position=4&search_layout=stack&type=product&tracking_id=3ecc166e-cc7f-401e-8926-b1671e3d6284"

# Row 8 - Sem Modelo (no politica value)
$ws.Range("C8").Value = "Sem Modelo"
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = "premium"
$ws.Range("G8").Value = "https://produto.mercadolivre.com.br/MLB-2772876015-filtro-anti-ruido-jfa-rca-eletromagnetico-stereo-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3ecc166e-cc7f-401e-8926-b1671e3d6284"

# Row 9 - Fonte 200A
$ws.Range("C9").Value = "FONTE 200A"
$ws.Range("D9").Value = "Igual"
$ws.Range("F9").Value = "premium"
$ws.Range("G9").Value = "https://produto.mercadolivre.com.br/MLB-3203737781-fonte-carregador-automotivo-jfa-storm-200-amperes-sci-bivolt-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3ecc166e-cc7f-401e-8926-b1671e3d6284"

# Row 10 - Fonte 90 Bob
$ws.Range("C10").Value = "FONTE 90 BOB"
$ws.Range("D10").Value = "Igual"
$ws.Range("F10").Value = "classico"
$ws.Range("G10").Value = "https://produto.mercadolivre.com.br/MLB-3240057378-fonte-carregador-jfa-90a-bob-storm-slim-bivolt-_JM#position%3D10%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3ecc166e-cc7f-401e-8926-b1671e3d6284"
